$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new FileError_Password / 2357 row first so the shared-string
# table picks up these values before the other changed/new strings
# (keeps shared-string ordering identical to the authored edit).
$ws.Range("A5").Value = "FileError_Password"
$ws.Range("B5").Value = "'2357"

# SourceFileError now points at the February file instead of January.
$ws.Range("B3").Value = "P:\96. Share Data\99. Other\13. IT\HOAI\QA_REPORT\NGUON\2.Bieu thong ke loi hang ngay 2.2023(New).xlsm"

# FileError_SheetName / value moved down; new value is "Main".
$ws.Range("A4").Value = "FileError_SheetName"
$ws.Range("B4").Value = "Main"

$ws.Range("B4").Select()
